$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 131; this shifts the existing rows 131..185
# down to 132..186 (Excel also copies formatting from the row above, as in
# the native Insert behaviour).
$ws.Rows("131").Insert()

# Populate the newly inserted row 131 with the new record.
$ws.Cells.Item(131, 1).Value = 4
$ws.Cells.Item(131, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(131, 3).Value = "Los Lagos"
$ws.Cells.Item(131, 4).Value = 44468
$ws.Cells.Item(131, 5).Value = 10
$ws.Cells.Item(131, 6).Value = 100114014
$ws.Cells.Item(131, 7).Value = "Betarraga"
$ws.Cells.Item(131, 8).Value = "Sin especificar"
$ws.Cells.Item(131, 9).Value = "Primera"
$ws.Cells.Item(131, 10).Value = 250
$ws.Cells.Item(131, 11).Value = 1200
$ws.Cells.Item(131, 12).Value = 1200
$ws.Cells.Item(131, 13).Value = 1200
$ws.Cells.Item(131, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(131, 15).Value = "Región del Maule"
$ws.Cells.Item(131, 16).Value = 240
$ws.Cells.Item(131, 17).Value = 5
$ws.Cells.Item(131, 18).Value = "Hortaliza"
